# Refresh the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# sheet with the latest scraped quote snapshot.
#
# Every cell in these two columns is stored as text in the workbook
# (prices keep thousands separators like "62.752.97", percentages keep
# their padding spaces like "  -1.45%  "), so a handful of the new price
# strings look like plain decimals (e.g. "565.98"). Those are written
# with a leading apostrophe so Excel keeps them as text instead of
# coercing them to numbers, then ClearFormats() drops the resulting
# "quote prefix" number format so the cell's style is left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.752.97'
$ws.Range('E2').Value = '  -1.45%  '
$ws.Range('D3').Value = '2.536.75'
$ws.Range('E3').Value = '  +2.79%  '
$ws.Range('E4').Value = '  +0.16%  '
$ws.Range('D5').Value = "'565.98"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.70%  '
$ws.Range('D6').Value = "'146.58"
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +2.42%  '
$ws.Range('E7').Value = '  +0.15%  '
$ws.Range('D8').Value = "'0.577"
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -2.40%  '
$ws.Range('D9').Value = '2.534.85'
$ws.Range('E9').Value = '  +2.71%  '
$ws.Range('E10').Value = '  -1.56%  '
$ws.Range('E11').Value = '  -2.87%  '
$ws.Range('E12').Value = '  +0.50%  '
$ws.Range('E13').Value = '  -0.89%  '
$ws.Range('D14').Value = "'26.90"
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +1.65%  '
$ws.Range('D15').Value = '2.992.95'
$ws.Range('E15').Value = '  +2.99%  '
$ws.Range('D16').Value = '62.768.19'
$ws.Range('E16').Value = '  -1.15%  '
$ws.Range('E17').Value = '  -2.07%  '
$ws.Range('D18').Value = '2.536.82'
$ws.Range('E18').Value = '  +2.97%  '
$ws.Range('E19').Value = '  +1.05%  '
$ws.Range('D20').Value = "'333.39"
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -2.86%  '
$ws.Range('E21').Value = '  -1.41%  '
$ws.Range('E22').Value = '  -1.14%  '
$ws.Range('E23').Value = '  -0.08%  '
$ws.Range('D24').Value = "'64.64"
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -1.82%  '
$ws.Range('E25').Value = '  -3.45%  '
$ws.Range('D26').Value = "'1.58"
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +2.57%  '
$ws.Range('D27').Value = "'0.998"
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -0.06%  '
$ws.Range('E28').Value = '  +10.56%  '
$ws.Range('E29').Value = '  +0.55%  '
$ws.Range('D30').Value = "'7.21"
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +4.61%  '
$ws.Range('D31').Value = '0.0₃0807'
$ws.Range('E31').Value = '  -1.76%  '
$ws.Range('E32').Value = '  -1.70%  '
$ws.Range('D33').Value = "'176.76"
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +0.91%  '
$ws.Range('E34').Value = '  +3.70%  '
$ws.Range('D35').Value = "'403.73"
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +8.13%  '
$ws.Range('D36').Value = "'18.86"
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -0.68%  '
$ws.Range('E37').Value = '  -1.60%  '
$ws.Range('E38').Value = '  -0.01%  '
$ws.Range('E39').Value = '  -3.97%  '
$ws.Range('E40').Value = '  -0.08%  '
$ws.Range('E41').Value = '  +0.18%  '
$ws.Range('E42').Value = '  -3.94%  '
$ws.Range('D43').Value = "'151.38"
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +0.69%  '
$ws.Range('E44').Value = '  -0.29%  '
$ws.Range('E45').Value = '  -1.34%  '
$ws.Range('E46').Value = '  -0.17%  '
$ws.Range('E47').Value = '  -1.35%  '
$ws.Range('D48').Value = "'0.0515"
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -2.14%  '
$ws.Range('E49').Value = '  +3.45%  '
$ws.Range('D50').Value = "'18.22"
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +0.12%  '
$ws.Range('E51').Value = '  +0.14%  '

Write-Output "done"
